# New Function: Add child/parent connection and update parameters for all widgets
#
# 1) Update the MP2153 (DCDC) voltageOutputRangeMin parameter (D7) from 1 to 0.5
# 2) Make the DCDC sheet the active sheet/tab, with D8 selected
#    (this also clears the previously-selected CONSUMER tab automatically)

$wb = $excel.ActiveWorkbook

$wsDCDC = $wb.Worksheets.Item("DCDC")

# Update the widget parameter value
$wsDCDC.Range("D7").Value = 0.5

# Make DCDC the active sheet and select D8, matching the saved view state
$wsDCDC.Activate()
$wsDCDC.Range("D8").Select() | Out-Null
